$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in A4
$ws.Range("A4").Value = 70

# Update selection to E8 (as seen in the diff)
$ws.Range("E8").Select()
